$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Data fixes (test method renames / typo fix) - per commit message:
# "refactored tests, moved test data preparation methods to BaseTest,
#  added BaseStep with enhanced logs"
# ---------------------------------------------------------------------

# Sheet "OwnersWithSameLastname": rename the search-pagination test case
# names used as fixture data (testCaseName column F).
$wsSameLastname = $wb.Worksheets.Item("OwnersWithSameLastname")

$wsSameLastname.Range("F30:F35").Value = "verifyPageNavigationByNumberTest"
$wsSameLastname.Range("F36:F41").Value = "verifyPageNavigationByArrowTest"

# Sheet "OwnerWithPets": fix the misspelled pet type "lizrd" -> "lizard".
$wsOwnerWithPets = $wb.Worksheets.Item("OwnerWithPets")
$wsOwnerWithPets.Range("H7").Value = "lizard"

# ---------------------------------------------------------------------
# Window / selection state: the workbook was last saved with the
# "OwnerWithPets" sheet active (instead of "Owner"), and each sheet's
# in-cell selection had moved around as the author navigated.
# ---------------------------------------------------------------------

$wsOwner = $wb.Worksheets.Item("Owner")
$wsOwner.Range("J16").Select()

$wsSameLastname.Activate()
$wsSameLastname.Range("A13").Select()
$wsSameLastname.Range("J39").Select()

$wsOwnerWithPets.Activate()
$wsOwnerWithPets.Range("H11").Select()
